$wb = $excel.ActiveWorkbook

# --- ALC sheet updates (hunk 0) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 500
$ws.Range("I19").Value = 500
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 500
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = -325
$ws.Range("N19").ClearContents()

# --- ALC sheet updates (hunk 1) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H42").Value = 633.25
$ws.Range("I42").Value = 686
$ws.Range("K42").Value = 2058
$ws.Range("M42").Value = -1828

# --- ALC sheet updates (hunk 2) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 3600
$ws.Range("J58").Value = 3500
$ws.Range("L58").Value = 10500
$ws.Range("N58").Value = -10800

# --- ALC sheet updates (hunk 3) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 944.8333
$ws.Range("I92").Value = 969.8333
$ws.Range("K92").Value = 969.8333
$ws.Range("M92").Value = 278.1667

# --- ALC sheet updates (hunk 4) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H111").Value = 782.5
$ws.Range("I111").Value = 782.5
$ws.Range("K111").Value = 2347.5
$ws.Range("M111").Value = 719.5

# --- ALC sheet updates (hunk 5) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 1302.375
$ws.Range("I132").Value = 1302.375
$ws.Range("K132").Value = 3907.125
$ws.Range("M132").Value = -1377.125

# --- ALC sheet updates (hunk 6) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 576.88
$ws.Range("I135").Value = 392.625
$ws.Range("K135").Value = 3533.625
$ws.Range("M135").Value = -998.625

# --- ALC sheet updates (hunk 7) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 2270.7144
$ws.Range("I137").Value = 1989
$ws.Range("J137").Value = 2975
$ws.Range("K137").Value = 5967
$ws.Range("L137").Value = 8925
$ws.Range("M137").Value = -3417
$ws.Range("N137").Value = -14025

# --- ALC sheet updates (hunk 8) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 4794.278
$ws.Range("I138").Value = 2869.8572
$ws.Range("K138").Value = 8609.571599999999
$ws.Range("M138").Value = -3469.571599999999

# --- ALC sheet updates (hunk 9) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 0
$ws.Range("I141").Value = 0
$ws.Range("K141").Value = 0
$ws.Range("M141").ClearContents()

# --- ARM sheet updates (hunk 10) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2029.6
$ws.Range("I2").Value = 2029.6
$ws.Range("K2").Value = 2029.6
$ws.Range("M2").Value = -1916.6

# --- ARM sheet updates (hunk 11) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2499.875
$ws.Range("I45").Value = 2519.8
$ws.Range("J45").Value = 2466.6667
$ws.Range("K45").Value = 2519.8
$ws.Range("L45").Value = 2466.6667
$ws.Range("M45").Value = -2142.8
$ws.Range("N45").Value = -3220.6667

# --- ARM sheet updates (hunk 12) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2292.3333
$ws.Range("I61").Value = 2261.8235
$ws.Range("K61").Value = 2261.8235
$ws.Range("M61").Value = -2049.8235

# --- ARM sheet updates (hunk 13) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 6286.2856
$ws.Range("I63").Value = 799.8
$ws.Range("K63").Value = 799.8
$ws.Range("M63").Value = -113.8

# --- ARM sheet updates (hunk 14) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 6286.2856
$ws.Range("I66").Value = 799.8
$ws.Range("K66").Value = 3999
$ws.Range("M66").Value = -567

# --- ARM sheet updates (hunk 15) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H95").Value = 24104
$ws.Range("J95").Value = 24104
$ws.Range("L95").Value = 24104
$ws.Range("N95").Value = -29596

# --- ARM sheet updates (hunk 16) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 3774.875
$ws.Range("I97").Value = 739.8
$ws.Range("J97").Value = 8833.333000000001
$ws.Range("K97").Value = 739.8
$ws.Range("L97").Value = 8833.333000000001
$ws.Range("M97").Value = -243.8
$ws.Range("N97").Value = -9825.333000000001

# --- ARM sheet updates (hunk 17) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 3131.6155
$ws.Range("I110").Value = 1335.5454
$ws.Range("J110").Value = 13010
$ws.Range("K110").Value = 1335.5454
$ws.Range("L110").Value = 13010
$ws.Range("M110").Value = 709.4546
$ws.Range("N110").Value = -17100

# --- ARM sheet updates (hunk 18) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 2029.6
$ws.Range("I116").Value = 2029.6
$ws.Range("K116").Value = 2029.6
$ws.Range("M116").Value = 264.4000000000001

# --- ARM sheet updates (hunk 19) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1622.6923
$ws.Range("I122").Value = 1341.25
$ws.Range("K122").Value = 4023.75
$ws.Range("M122").Value = -1573.75

# --- ARM sheet updates (hunk 20) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 2292.3333
$ws.Range("I136").Value = 2261.8235
$ws.Range("K136").Value = 6785.470499999999
$ws.Range("M136").Value = -4235.470499999999

# --- BSM sheet updates (hunk 21) ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2029.6
$ws.Range("I3").Value = 2029.6
$ws.Range("K3").Value = 2029.6
$ws.Range("M3").Value = -1915.6

# --- BSM sheet updates (hunk 22) ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1662.2727
$ws.Range("I107").Value = 1616.3334
$ws.Range("J107").Value = 1717.4
$ws.Range("K107").Value = 1616.3334
$ws.Range("L107").Value = 1717.4
$ws.Range("M107").Value = 303.6666
$ws.Range("N107").Value = -5557.4

# --- CRP sheet updates (hunk 23) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3890.111
$ws.Range("I58").Value = 2301.8333
$ws.Range("K58").Value = 2301.8333
$ws.Range("M58").Value = -2098.8333

# --- CRP sheet updates (hunk 24) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H60").Value = 47732.5
$ws.Range("I60").Value = 0
$ws.Range("J60").Value = 47732.5
$ws.Range("K60").Value = 0
$ws.Range("L60").Value = 47732.5
$ws.Range("M60").ClearContents()
$ws.Range("N60").Value = -48754.5

# --- CRP sheet updates (hunk 25) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1981.125
$ws.Range("I134").Value = 999.25
$ws.Range("K134").Value = 2997.75
$ws.Range("M134").Value = -462.75

# --- CRP sheet updates (hunk 26) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 3890.111
$ws.Range("I136").Value = 2301.8333
$ws.Range("K136").Value = 6905.499899999999
$ws.Range("M136").Value = -4355.499899999999

# --- CUL sheet updates (hunk 27) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 13400
$ws.Range("I80").Value = 8000
$ws.Range("J80").Value = 14000
$ws.Range("K80").Value = 24000
$ws.Range("L80").Value = 42000
$ws.Range("M80").Value = -23064
$ws.Range("N80").Value = -43872

# --- CUL sheet updates (hunk 28) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H83").Value = 13400
$ws.Range("I83").Value = 8000
$ws.Range("J83").Value = 14000
$ws.Range("K83").Value = 72000
$ws.Range("L83").Value = 126000
$ws.Range("M83").Value = -67320
$ws.Range("N83").Value = -135360

# --- CUL sheet updates (hunk 29) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H109").Value = 4201.6665
$ws.Range("I109").Value = 1973.3334
$ws.Range("K109").Value = 5920.0002
$ws.Range("M109").Value = -4880.0002

# --- GSM sheet updates (hunk 30) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H17").Value = 0
$ws.Range("J17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("N17").ClearContents()

# --- LTW sheet updates (hunk 31) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3282.375
$ws.Range("I7").Value = 3282.375
$ws.Range("K7").Value = 3282.375
$ws.Range("M7").Value = -3170.375

# --- LTW sheet updates (hunk 32) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 8219.200000000001
$ws.Range("J22").Value = 8685.714
$ws.Range("L22").Value = 8685.714
$ws.Range("N22").Value = -9275.714

# --- LTW sheet updates (hunk 33) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 8219.200000000001
$ws.Range("J27").Value = 8685.714
$ws.Range("L27").Value = 8685.714
$ws.Range("N27").Value = -8899.714

# --- LTW sheet updates (hunk 34) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2183.625
$ws.Range("I68").Value = 2044.8334
$ws.Range("K68").Value = 2044.8334
$ws.Range("M68").Value = -1295.8334

# --- LTW sheet updates (hunk 35) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 2183.625
$ws.Range("I71").Value = 2044.8334
$ws.Range("K71").Value = 10224.167
$ws.Range("M71").Value = -6480.166999999999

# --- LTW sheet updates (hunk 36) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3861.75
$ws.Range("I122").Value = 3861.75
$ws.Range("K122").Value = 11585.25
$ws.Range("M122").Value = -9135.25

# --- LTW sheet updates (hunk 37) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 3282.375
$ws.Range("I126").Value = 3282.375
$ws.Range("K126").Value = 9847.125
$ws.Range("M126").Value = -7377.125

# --- WVR sheet updates (hunk 38) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 28880
$ws.Range("J54").Value = 33506.668
$ws.Range("L54").Value = 33506.668
$ws.Range("N54").Value = -34546.668

# --- WVR sheet updates (hunk 39) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 7194.3
$ws.Range("I62").Value = 4133.3335
$ws.Range("J62").Value = 8506.143
$ws.Range("K62").Value = 4133.3335
$ws.Range("L62").Value = 8506.143
$ws.Range("M62").Value = -3509.3335
$ws.Range("N62").Value = -9754.143

# --- WVR sheet updates (hunk 40) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 7194.3
$ws.Range("I65").Value = 4133.3335
$ws.Range("J65").Value = 8506.143
$ws.Range("K65").Value = 20666.6675
$ws.Range("L65").Value = 42530.715
$ws.Range("M65").Value = -17546.6675
$ws.Range("N65").Value = -48770.715

# --- WVR sheet updates (hunk 41) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 792.5833
$ws.Range("I100").Value = 748.7143
$ws.Range("J100").Value = 854
$ws.Range("K100").Value = 1497.4286
$ws.Range("L100").Value = 1708
$ws.Range("M100").Value = -956.4286
$ws.Range("N100").Value = -2790

# --- WVR sheet updates (hunk 42) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1101.1428
$ws.Range("I122").Value = 1118.1666
$ws.Range("J122").Value = 999
$ws.Range("K122").Value = 3354.4998
$ws.Range("L122").Value = 2997
$ws.Range("M122").Value = -904.4998000000001
$ws.Range("N122").Value = -7897

# --- WVR sheet updates (hunk 43) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()

# --- WVR sheet updates (hunk 44) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1686.2727
$ws.Range("I126").Value = 1601
$ws.Range("K126").Value = 4803
$ws.Range("M126").Value = -2333

# --- WVR sheet updates (hunk 45) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3018.7368
$ws.Range("I132").Value = 2500.6
$ws.Range("K132").Value = 7501.799999999999
$ws.Range("M132").Value = -4971.799999999999

# --- WVR sheet updates (hunk 46) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 8495.6875
$ws.Range("I136").Value = 7764
$ws.Range("J136").Value = 11666.333
$ws.Range("K136").Value = 23292
$ws.Range("L136").Value = 34998.999
$ws.Range("M136").Value = -20742
$ws.Range("N136").Value = -40098.999
